$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 225 - this shifts existing rows 225-261 down to 226-262
# and copies formatting (e.g. date number format) from the row above.
$ws.Rows.Item(225).Insert()

# Populate the newly inserted row 225 with the new record's data.
$ws.Range("A225").Value2 = 10
$ws.Range("B225").Value2 = "Vega Modelo de Temuco"
$ws.Range("C225").Value2 = "La Araucanía"
$ws.Range("D225").Value2 = 45258
$ws.Range("E225").Value2 = 9
$ws.Range("F225").Value2 = "Fruta"
$ws.Range("G225").Value2 = 100107
$ws.Range("H225").Value2 = "Otros"
$ws.Range("I225").Value2 = 100107002
$ws.Range("J225").Value2 = "Chirimoya"
$ws.Range("K225").Value2 = "Cultivar IV Región"
$ws.Range("L225").Value2 = "Primera"
$ws.Range("M225").Value2 = 75
$ws.Range("N225").Value2 = 2400
$ws.Range("O225").Value2 = 2400
$ws.Range("P225").Value2 = 2400
$ws.Range("Q225").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R225").Value2 = "Provincia del Elquí"
$ws.Range("S225").Value2 = 2400
$ws.Range("T225").Value2 = 1
